# Updated cryptos list (prices / 1h volume change) - GitHub Actions refresh.
# All Price/Volume cells in this sheet are stored as text (inlineStr), even
# though many look like plain numbers (e.g. "0.9985", "5.010", "243.29").
# Writing such a string straight into Range.Value would make Excel parse it
# as a real number and silently reformat it (e.g. "5.010" -> 5.01), which
# would not match the source data. To keep them as literal text without
# leaving any stray cell-style (NumberFormat) behind, we temporarily force
# the cell to Text format, assign the value, then restore the "Normal"
# style so the saved XML has no extra style attribute on the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextValue 'D2' '29.448.17'
Set-TextValue 'E2' '  +0.90%  '
Set-TextValue 'D3' '1.839.07'
Set-TextValue 'E3' '  +0.04%  '
Set-TextValue 'D4' '0.9985'
Set-TextValue 'E4' '  -1.11%  '
Set-TextValue 'D5' '243.29'
Set-TextValue 'D6' '0.6274'
Set-TextValue 'E6' '  +1.58%  '
Set-TextValue 'D7' '0.9993'
Set-TextValue 'E7' '  -1.10%  '
Set-TextValue 'D8' '0.07404'
Set-TextValue 'E8' '  -0.46%  '
Set-TextValue 'D9' '0.2957'
Set-TextValue 'E9' '  +0.50%  '
Set-TextValue 'D10' '23.49'
Set-TextValue 'E10' '  +2.46%  '
Set-TextValue 'D11' '0.07646'
Set-TextValue 'E11' '  -0.86%  '
Set-TextValue 'D12' '1.835.10'
Set-TextValue 'E12' '  -0.02%  '
Set-TextValue 'D13' '5.010'
Set-TextValue 'E13' '  +0.59%  '
Set-TextValue 'D14' '0.6757'
Set-TextValue 'E14' '  +0.74%  '
Set-TextValue 'D15' '83.53'
Set-TextValue 'E15' '  +0.97%  '
Set-TextValue 'D16' '0.000009259'
Set-TextValue 'E16' '  +1.96%  '
Set-TextValue 'D17' '5.903'
Set-TextValue 'E17' '  +0.55%  '
Set-TextValue 'D18' '29.400.03'
Set-TextValue 'E18' '  +0.68%  '
Set-TextValue 'D19' '2.082.82'
Set-TextValue 'E19' '  +0.19%  '
Set-TextValue 'D20' '236.88'
Set-TextValue 'E20' '  +0.11%  '
Set-TextValue 'E21' '  -0.33%  '
Set-TextValue 'D22' '0.9993'
Set-TextValue 'E22' '  -1.13%  '
Set-TextValue 'D23' '7.322'
Set-TextValue 'E23' '  +2.41%  '
Set-TextValue 'D24' '0.9997'
Set-TextValue 'E24' '  -1.32%  '
Set-TextValue 'D25' '158.90'
Set-TextValue 'E25' '  -0.38%  '
Set-TextValue 'D26' '0.1413'
Set-TextValue 'E26' '  -0.33%  '
Set-TextValue 'D27' '8.498'
Set-TextValue 'E27' '  +0.02%  '
Set-TextValue 'D28' '17.72'
Set-TextValue 'E28' '  -0.68%  '
Set-TextValue 'B29' 'Hedera'
Set-TextValue 'C29' 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue 'D29' '0.06047'
Set-TextValue 'E29' '  +8.87%  '
Set-TextValue 'B30' 'PancakeSwap'
Set-TextValue 'C30' 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
Set-TextValue 'D30' '1.493'
Set-TextValue 'E30' '  -0.74%  '
Set-TextValue 'D31' '1.236'
Set-TextValue 'E31' '  +1.58%  '
Set-TextValue 'D32' '4.096'
Set-TextValue 'E32' '  -0.36%  '
Set-TextValue 'D33' '4.106'
Set-TextValue 'E33' '  -0.59%  '
Set-TextValue 'D34' '1.872'
Set-TextValue 'E34' '  +1.04%  '
Set-TextValue 'E35' '  +0.10%  '
Set-TextValue 'D36' '0.7251'
Set-TextValue 'E36' '  -2.03%  '
Set-TextValue 'D37' '2.611'
Set-TextValue 'E37' '  -1.72%  '
Set-TextValue 'D38' '2.873'
Set-TextValue 'E38' '  +1.57%  '
Set-TextValue 'D39' '1.216.26'
Set-TextValue 'E39' '  +0.98%  '
Set-TextValue 'E40' '  -0.83%  '
Set-TextValue 'D41' '6.256'
Set-TextValue 'E41' '  -2.59%  '
Set-TextValue 'D42' '0.9109'
Set-TextValue 'E42' '  +0.92%  '
Set-TextValue 'E43' '  -0.94%  '
Set-TextValue 'D44' '1.994.99'
Set-TextValue 'E44' '  +0.57%  '
Set-TextValue 'D45' '102.02'
Set-TextValue 'E45' '  +0.67%  '
Set-TextValue 'D46' '65.42'
Set-TextValue 'E46' '  +0.76%  '
Set-TextValue 'D47' '0.00000000121'
Set-TextValue 'E47' '  -1.60%  '
Set-TextValue 'D48' '0.5068'
Set-TextValue 'E48' '  -1.27%  '
Set-TextValue 'D49' '9.246'
Set-TextValue 'E49' '  +1.82%  '
Set-TextValue 'D50' '0.4060'
Set-TextValue 'E50' '  +0.72%  '
Set-TextValue 'D51' '0.1138'
Set-TextValue 'E51' '  +3.03%  '
